# Apply the edit described in the diff:
#  - Row 9  (Sarah Swenson): clear title/office/phone/email (F9:I9)
#  - Row 12 (Brett Wilder): fill in title/office/phone/email (F12:I12) and degrees (L12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: clear out the previously-filled columns F-I
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = ""

# Row 12: populate title, office, phone, email, and degrees
$ws.Range("F12").Value = "Assistant Professor & Area Extension Educator — Farm Business Management"
$ws.Range("G12").Value = "Room C1"
$ws.Range("H12").Value = "208-885-0263"
$ws.Range("I12").Value = "bwilder@uidaho.edu"
$ws.Range("L12").Value = "['M.S., University of Idaho, 2019', 'B.S., University of Idaho, 2017']"
